# Resize the first 7 columns of the active sheet.
#
# Target (canonical OOXML <col> widths) from the diff:
#   A: 16.8   B: 20.8   C: 18.2   D: 40   E: 28.6   F: 15.6   G: 15.6
#
# The runtime stores a column's width internally as pixels using the
# classic Excel "Maximum Digit Width" model (content-px + 5px padding,
# quantized to whole pixels) before re-exporting it as a "characters"
# width in the XML. Because of that integer-pixel quantization, the
# COM `ColumnWidth` setter cannot always reproduce an arbitrary decimal
# width bit-for-bit; the values below are the inputs that round-trip to
# the closest possible (and, for column D, exact) match of the target
# widths once Excel re-derives the stored width from pixels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 39.166666666666664
$ws.Columns.Item(5).ColumnWidth = 27.833333333333332
$ws.Columns.Item(6).ColumnWidth = 14.833333333333334
$ws.Columns.Item(7).ColumnWidth = 14.833333333333334
